$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.337.04"
$ws.Range("E2").Value = "  +0.24%  "
$ws.Range("D3").Value = "1.611.90"
$ws.Range("E3").Value = "  +0.26%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'213.11"
$ws.Range("E5").Value = "  -0.06%  "
$ws.Range("E6").Value = "  -0.06%  "
$ws.Range("D7").Value = "'0.486"
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  -0.43%  "
$ws.Range("D10").Value = "'18.51"
$ws.Range("E10").Value = "  +2.31%  "
$ws.Range("D11").Value = "'0.0813"
$ws.Range("E11").Value = "  -1.31%  "
$ws.Range("D12").Value = "1.838.22"
$ws.Range("D13").Value = "1.607.12"
$ws.Range("E13").Value = "  +0.00%  "
$ws.Range("D14").Value = "'4.04"
$ws.Range("E14").Value = "  +0.63%  "
$ws.Range("D15").Value = "'0.516"
$ws.Range("E15").Value = "  +0.79%  "
$ws.Range("D16").Value = "26.345.24"
$ws.Range("E16").Value = "  +0.41%  "
$ws.Range("D17").Value = "'61.82"
$ws.Range("E17").Value = "  +1.87%  "
$ws.Range("D18").Value = "0.0₃0730"
$ws.Range("E18").Value = "  +0.66%  "
$ws.Range("E19").Value = "  -0.12%  "
$ws.Range("D20").Value = "'203.39"
$ws.Range("E20").Value = "  +2.29%  "
$ws.Range("E21").Value = "  +1.01%  "
$ws.Range("E22").Value = "  -0.09%  "
$ws.Range("E23").Value = "  +0.51%  "
$ws.Range("E24").Value = "  +7.94%  "
$ws.Range("D25").Value = "'144.42"
$ws.Range("E25").Value = "  +1.11%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("E27").Value = "  -2.93%  "
$ws.Range("E28").Value = "  +0.35%  "
$ws.Range("D29").Value = "'6.59"
$ws.Range("E29").Value = "  +1.88%  "
$ws.Range("D30").Value = "'0.0490"
$ws.Range("E30").Value = "  +3.99%  "
$ws.Range("E31").Value = "  +0.14%  "
$ws.Range("E32").Value = "  +1.78%  "
$ws.Range("E33").Value = "  -1.53%  "
$ws.Range("D34").Value = "'2.44"
$ws.Range("E34").Value = "  +3.35%  "
$ws.Range("D36").Value = "1.162.51"
$ws.Range("E36").Value = "  +5.12%  "
$ws.Range("D37").Value = "'0.0166"
$ws.Range("E37").Value = "  +9.40%  "
$ws.Range("E38").Value = "  -0.07%  "
$ws.Range("D39").Value = "'0.795"
$ws.Range("E39").Value = "  +0.80%  "
$ws.Range("E40").Value = "  -0.28%  "
$ws.Range("D41").Value = "'0.501"
$ws.Range("E41").Value = "  +0.29%  "
$ws.Range("D42").Value = "'0.786"
$ws.Range("E42").Value = "  +1.44%  "
$ws.Range("D43").Value = "'5.25"
$ws.Range("E43").Value = "  +2.91%  "
$ws.Range("D44").Value = "1.751.52"
$ws.Range("E44").Value = "  +0.44%  "
$ws.Range("D45").Value = "'91.76"
$ws.Range("E45").Value = "  -1.15%  "
$ws.Range("E46").Value = "  -0.27%  "
$ws.Range("D47").Value = "'54.47"
$ws.Range("E47").Value = "  +1.68%  "
$ws.Range("E48").Value = "  +0.17%  "
$ws.Range("D49").Value = "0.0₇0979"
$ws.Range("E49").Value = "  -13.84%  "
$ws.Range("E50").Value = "  -0.64%  "
$ws.Range("E51").Value = "  -0.11%  "
